$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "M1"
$ws.Cells.Item(2, 2).Value = "Il1f5"
$ws.Cells.Item(2, 3).Value = "Il1rl2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.02506066666666667
$ws.Cells.Item(2, 8).Value = 0.075182
$ws.Cells.Item(2, 9).Value = 0.005556484732231195
$ws.Cells.Item(2, 10).Value = 0.005556484732231196
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.572473333333333
$ws.Cells.Item(2, 14).Value = 4.71742
$ws.Cells.Item(2, 15).Value = 0.0623955221013184
$ws.Cells.Item(2, 16).Value = 0.06239552210131842
$ws.Cells.Item(2, 17).Value = 0.03940723004888888
$ws.Cells.Item(2, 18).Value = 0.35466507044
$ws.Cells.Item(2, 19).Value = 0.0003466997659155698
$ws.Cells.Item(2, 20).Value = 0.00034669976591557

# Row 3
$ws.Cells.Item(3, 1).Value = "M1"
$ws.Cells.Item(3, 2).Value = "Il1f5"
$ws.Cells.Item(3, 3).Value = "Il1rl2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.02506066666666667
$ws.Cells.Item(3, 8).Value = 0.075182
$ws.Cells.Item(3, 9).Value = 0.005556484732231195
$ws.Cells.Item(3, 10).Value = 0.005556484732231196
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 16.615168
$ws.Cells.Item(3, 14).Value = 49.84550400000001
$ws.Cells.Item(3, 15).Value = 0.6592875441413644
$ws.Cells.Item(3, 16).Value = 0.6592875441413647
$ws.Cells.Item(3, 17).Value = 0.4163871868586667
$ws.Cells.Item(3, 18).Value = 3.747484681728
$ws.Cells.Item(3, 19).Value = 0.003663321173171692
$ws.Cells.Item(3, 20).Value = 0.003663321173171694

# Row 4
$ws.Cells.Item(4, 1).Value = "M1"
$ws.Cells.Item(4, 2).Value = "Il1f5"
$ws.Cells.Item(4, 3).Value = "Il1rl2"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.02506066666666667
$ws.Cells.Item(4, 8).Value = 0.075182
$ws.Cells.Item(4, 9).Value = 0.005556484732231195
$ws.Cells.Item(4, 10).Value = 0.005556484732231196
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.354957666666667
$ws.Cells.Item(4, 14).Value = 4.064873
$ws.Cells.Item(4, 15).Value = 0.05376453084748706
$ws.Cells.Item(4, 16).Value = 0.05376453084748708
$ws.Cells.Item(4, 17).Value = 0.03395614243177777
$ws.Cells.Item(4, 18).Value = 0.3056052818859999
$ws.Cells.Item(4, 19).Value = 0.000298741794789635
$ws.Cells.Item(4, 20).Value = 0.0002987417947896351

# Row 5
$ws.Cells.Item(5, 1).Value = "M1"
$ws.Cells.Item(5, 2).Value = "Il1f5"
$ws.Cells.Item(5, 3).Value = "Il1rl2"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.02506066666666667
$ws.Cells.Item(5, 8).Value = 0.075182
$ws.Cells.Item(5, 9).Value = 0.005556484732231195
$ws.Cells.Item(5, 10).Value = 0.005556484732231196
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.87712
$ws.Cells.Item(5, 14).Value = 5.63136
$ws.Cells.Item(5, 15).Value = 0.07448385925791649
$ws.Cells.Item(5, 16).Value = 0.07448385925791652
$ws.Cells.Item(5, 17).Value = 0.04704187861333333
$ws.Cells.Item(5, 18).Value = 0.42337690752
$ws.Cells.Item(5, 19).Value = 0.0004138684267642702
$ws.Cells.Item(5, 20).Value = 0.0004138684267642704

# Row 6
$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 2).Value = "Il1f5"
$ws.Cells.Item(6, 3).Value = "Il1rl2"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.02506066666666667
$ws.Cells.Item(6, 8).Value = 0.075182
$ws.Cells.Item(6, 9).Value = 0.005556484732231195
$ws.Cells.Item(6, 10).Value = 0.005556484732231196
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.367415333333334
$ws.Cells.Item(6, 14).Value = 7.102246000000001
$ws.Cells.Item(6, 15).Value = 0.09393870956200642
$ws.Cells.Item(6, 16).Value = 0.09393870956200645
$ws.Cells.Item(6, 17).Value = 0.05932900653022222
$ws.Cells.Item(6, 18).Value = 0.5339610587720001
$ws.Cells.Item(6, 19).Value = 0.0005219690054467893
$ws.Cells.Item(6, 20).Value = 0.0005219690054467896

# Row 7
$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Il1f5"
$ws.Cells.Item(7, 3).Value = "Il1rl2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.02506066666666667
$ws.Cells.Item(7, 8).Value = 0.075182
$ws.Cells.Item(7, 9).Value = 0.005556484732231195
$ws.Cells.Item(7, 10).Value = 0.005556484732231196
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.414567333333333
$ws.Cells.Item(7, 14).Value = 4.243702
$ws.Cells.Item(7, 15).Value = 0.05612983408990701
$ws.Cells.Item(7, 16).Value = 0.05612983408990703
$ws.Cells.Item(7, 17).Value = 0.03545000041822222
$ws.Cells.Item(7, 18).Value = 0.319050003764
$ws.Cells.Item(7, 19).Value = 0.0003118845661432384
$ws.Cells.Item(7, 20).Value = 0.0003118845661432385

# Row 8
$ws.Cells.Item(8, 1).Value = "Neutro"
$ws.Cells.Item(8, 2).Value = "Il1f5"
$ws.Cells.Item(8, 3).Value = "Il1rl2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.485105
$ws.Cells.Item(8, 8).Value = 13.455315
$ws.Cells.Item(8, 9).Value = 0.9944435152677689
$ws.Cells.Item(8, 10).Value = 0.9944435152677689
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.572473333333333
$ws.Cells.Item(8, 14).Value = 4.71742
$ws.Cells.Item(8, 15).Value = 0.0623955221013184
$ws.Cells.Item(8, 16).Value = 0.06239552210131842
$ws.Cells.Item(8, 17).Value = 7.052708009699999
$ws.Cells.Item(8, 18).Value = 63.47437208729999
$ws.Cells.Item(8, 19).Value = 0.06204882233540283
$ws.Cells.Item(8, 20).Value = 0.06204882233540285

# Row 9
$ws.Cells.Item(9, 1).Value = "Neutro"
$ws.Cells.Item(9, 2).Value = "Il1f5"
$ws.Cells.Item(9, 3).Value = "Il1rl2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.485105
$ws.Cells.Item(9, 8).Value = 13.455315
$ws.Cells.Item(9, 9).Value = 0.9944435152677689
$ws.Cells.Item(9, 10).Value = 0.9944435152677689
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 16.615168
$ws.Cells.Item(9, 14).Value = 49.84550400000001
$ws.Cells.Item(9, 15).Value = 0.6592875441413644
$ws.Cells.Item(9, 16).Value = 0.6592875441413647
$ws.Cells.Item(9, 17).Value = 74.52077307264
$ws.Cells.Item(9, 18).Value = 670.68695765376
$ws.Cells.Item(9, 19).Value = 0.6556242229681927
$ws.Cells.Item(9, 20).Value = 0.655624222968193

# Row 10
$ws.Cells.Item(10, 1).Value = "Neutro"
$ws.Cells.Item(10, 2).Value = "Il1f5"
$ws.Cells.Item(10, 3).Value = "Il1rl2"
$ws.Cells.Item(10, 4).Value = "M1"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.485105
$ws.Cells.Item(10, 8).Value = 13.455315
$ws.Cells.Item(10, 9).Value = 0.9944435152677689
$ws.Cells.Item(10, 10).Value = 0.9944435152677689
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.354957666666667
$ws.Cells.Item(10, 14).Value = 4.064873
$ws.Cells.Item(10, 15).Value = 0.05376453084748706
$ws.Cells.Item(10, 16).Value = 0.05376453084748708
$ws.Cells.Item(10, 17).Value = 6.077127405554999
$ws.Cells.Item(10, 18).Value = 54.69414664999499
$ws.Cells.Item(10, 19).Value = 0.05346578905269743
$ws.Cells.Item(10, 20).Value = 0.05346578905269744

# Row 11
$ws.Cells.Item(11, 1).Value = "Neutro"
$ws.Cells.Item(11, 2).Value = "Il1f5"
$ws.Cells.Item(11, 3).Value = "Il1rl2"
$ws.Cells.Item(11, 4).Value = "M2"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.485105
$ws.Cells.Item(11, 8).Value = 13.455315
$ws.Cells.Item(11, 9).Value = 0.9944435152677689
$ws.Cells.Item(11, 10).Value = 0.9944435152677689
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.87712
$ws.Cells.Item(11, 14).Value = 5.63136
$ws.Cells.Item(11, 15).Value = 0.07448385925791649
$ws.Cells.Item(11, 16).Value = 0.07448385925791652
$ws.Cells.Item(11, 17).Value = 8.419080297599999
$ws.Cells.Item(11, 18).Value = 75.7717226784
$ws.Cells.Item(11, 19).Value = 0.07406999083115223
$ws.Cells.Item(11, 20).Value = 0.07406999083115226

# Row 12
$ws.Cells.Item(12, 1).Value = "Neutro"
$ws.Cells.Item(12, 2).Value = "Il1f5"
$ws.Cells.Item(12, 3).Value = "Il1rl2"
$ws.Cells.Item(12, 4).Value = "Neutro"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.485105
$ws.Cells.Item(12, 8).Value = 13.455315
$ws.Cells.Item(12, 9).Value = 0.9944435152677689
$ws.Cells.Item(12, 10).Value = 0.9944435152677689
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.367415333333334
$ws.Cells.Item(12, 14).Value = 7.102246000000001
$ws.Cells.Item(12, 15).Value = 0.09393870956200642
$ws.Cells.Item(12, 16).Value = 0.09393870956200645
$ws.Cells.Item(12, 17).Value = 10.61810634861
$ws.Cells.Item(12, 18).Value = 95.56295713749
$ws.Cells.Item(12, 19).Value = 0.09341674055655964
$ws.Cells.Item(12, 20).Value = 0.09341674055655967

# Row 13
$ws.Cells.Item(13, 1).Value = "Neutro"
$ws.Cells.Item(13, 2).Value = "Il1f5"
$ws.Cells.Item(13, 3).Value = "Il1rl2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.485105
$ws.Cells.Item(13, 8).Value = 13.455315
$ws.Cells.Item(13, 9).Value = 0.9944435152677689
$ws.Cells.Item(13, 10).Value = 0.9944435152677689
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 1.414567333333333
$ws.Cells.Item(13, 14).Value = 4.243702
$ws.Cells.Item(13, 15).Value = 0.05612983408990701
$ws.Cells.Item(13, 16).Value = 0.05612983408990703
$ws.Cells.Item(13, 17).Value = 6.344483019569999
$ws.Cells.Item(13, 18).Value = 57.10034717612999
$ws.Cells.Item(13, 19).Value = 0.05581794952376377
$ws.Cells.Item(13, 20).Value = 0.05581794952376379

